$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Tiny precision refresh on an existing value (recompute artifact) ---
$ws.Range("B2").Value = -326.15795207023621

# --- Add a new "mean" summary row (row 9) under the data table ---

# Bring over the bold/bordered header look for the new row label cell (A9)
$ws.Range("A1").Copy()
$ws.Range("A9").PasteSpecial(-4122)

# Bring over the numeric formatting used throughout the data rows (B9:F9)
$ws.Range("B2:F2").Copy()
$ws.Range("B9:F9").PasteSpecial(-4122)

# Row label
$ws.Range("A9").Value = "mean"

# Column averages - B9 seeds the formula, C9:F9 fill in as a shared formula
$ws.Range("B9").Formula = "=AVERAGE(B2:B8)"
$ws.Range("C9:F9").Formula = "=AVERAGE(C2:C8)"

# Clear the clipboard marching ants left over from the copy/paste operations
$excel.CutCopyMode = 0

# Reflect the newly extended data range in the current selection
$ws.Range("B2:F9").Select()

Write-Output "done"
